# This script updates the "想去人数" (interested-count) figures in column F
# across the four worksheets of the workbook, reflecting a newer scrape of
# the source data. Only the numeric values in column F change; everything
# else (labels, other columns, formatting) stays the same.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F10").Value = 5596
$ws.Range("F18").Value = 4745
$ws.Range("F21").Value = 2385
$ws.Range("F24").Value = 1173
$ws.Range("F31").Value = 2006
$ws.Range("F37").Value = 595
$ws.Range("F40").Value = 186
$ws.Range("F41").Value = 1638

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 52
$ws.Range("F20").Value = 138
$ws.Range("F23").Value = 142

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1331
$ws.Range("F10").Value = 1761
$ws.Range("F11").Value = 2243
$ws.Range("F12").Value = 685
$ws.Range("F13").Value = 561

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1331
$ws.Range("F8").Value = 2243
$ws.Range("F9").Value = 5596
$ws.Range("F10").Value = 685
$ws.Range("F21").Value = 4745
$ws.Range("F22").Value = 2385
$ws.Range("F23").Value = 1173
$ws.Range("F29").Value = 52
$ws.Range("F33").Value = 2006
$ws.Range("F35").Value = 138
$ws.Range("F38").Value = 142
$ws.Range("F39").Value = 595
$ws.Range("F43").Value = 1638
